# Update Long27_DataComp.xlsx cells per the source diff:
# - Column C/D (M2_Len / FX_Len) counters incremented by 1 on affected rows.
# - Columns E-H (M2_1stDate, M2_LastDate, FX_1stDate, FX_LastDate) rolled
#   forward by one month (date serials) on affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = "G2"; Value = 30074 },
    @{ Cell = "H2"; Value = 45261 },
    @{ Cell = "E3"; Value = 30011 },
    @{ Cell = "F3"; Value = 45200 },
    @{ Cell = "E4"; Value = 30011 },
    @{ Cell = "F4"; Value = 45200 },
    @{ Cell = "G4"; Value = 30074 },
    @{ Cell = "H4"; Value = 45261 },
    @{ Cell = "G5"; Value = 30074 },
    @{ Cell = "H5"; Value = 45261 },
    @{ Cell = "C6"; Value = 443 },
    @{ Cell = "F6"; Value = 45200 },
    @{ Cell = "G6"; Value = 30074 },
    @{ Cell = "H6"; Value = 45261 },
    @{ Cell = "G7"; Value = 30074 },
    @{ Cell = "H7"; Value = 45261 },
    @{ Cell = "D8"; Value = 410 },
    @{ Cell = "E8"; Value = 30011 },
    @{ Cell = "F8"; Value = 45200 },
    @{ Cell = "H8"; Value = 45261 },
    @{ Cell = "E9"; Value = 30011 },
    @{ Cell = "F9"; Value = 45200 },
    @{ Cell = "G9"; Value = 30074 },
    @{ Cell = "H9"; Value = 45261 },
    @{ Cell = "D10"; Value = 483 },
    @{ Cell = "E10"; Value = 30011 },
    @{ Cell = "F10"; Value = 45200 },
    @{ Cell = "H10"; Value = 45261 },
    @{ Cell = "G11"; Value = 30074 },
    @{ Cell = "H11"; Value = 45261 },
    @{ Cell = "C12"; Value = 371 },
    @{ Cell = "D12"; Value = 352 },
    @{ Cell = "F12"; Value = 45200 },
    @{ Cell = "H12"; Value = 45261 },
    @{ Cell = "C13"; Value = 467 },
    @{ Cell = "F13"; Value = 45200 },
    @{ Cell = "G13"; Value = 30074 },
    @{ Cell = "H13"; Value = 45261 },
    @{ Cell = "D14"; Value = 396 },
    @{ Cell = "H14"; Value = 45261 },
    @{ Cell = "G15"; Value = 30042 },
    @{ Cell = "H15"; Value = 45261 },
    @{ Cell = "C16"; Value = 455 },
    @{ Cell = "D16"; Value = 410 },
    @{ Cell = "F16"; Value = 45200 },
    @{ Cell = "H16"; Value = 45261 },
    @{ Cell = "C17"; Value = 370 },
    @{ Cell = "D17"; Value = 394 },
    @{ Cell = "F17"; Value = 45200 },
    @{ Cell = "H17"; Value = 45261 },
    @{ Cell = "E18"; Value = 30011 },
    @{ Cell = "F18"; Value = 45200 },
    @{ Cell = "G18"; Value = 30074 },
    @{ Cell = "H18"; Value = 45261 },
    @{ Cell = "D19"; Value = 398 },
    @{ Cell = "E19"; Value = 30011 },
    @{ Cell = "F19"; Value = 45200 },
    @{ Cell = "H19"; Value = 45261 },
    @{ Cell = "C20"; Value = 481 },
    @{ Cell = "F20"; Value = 45200 },
    @{ Cell = "G20"; Value = 30074 },
    @{ Cell = "H20"; Value = 45261 },
    @{ Cell = "E21"; Value = 30011 },
    @{ Cell = "F21"; Value = 45200 },
    @{ Cell = "G21"; Value = 30074 },
    @{ Cell = "H21"; Value = 45261 },
    @{ Cell = "D22"; Value = 380 },
    @{ Cell = "H22"; Value = 45261 },
    @{ Cell = "D23"; Value = 217 },
    @{ Cell = "E23"; Value = 30011 },
    @{ Cell = "F23"; Value = 45200 },
    @{ Cell = "H23"; Value = 45261 },
    @{ Cell = "C24"; Value = 394 },
    @{ Cell = "F24"; Value = 45200 },
    @{ Cell = "G24"; Value = 30074 },
    @{ Cell = "H24"; Value = 45261 },
    @{ Cell = "E25"; Value = 30011 },
    @{ Cell = "F25"; Value = 45200 },
    @{ Cell = "G25"; Value = 30074 },
    @{ Cell = "H25"; Value = 45261 },
    @{ Cell = "D26"; Value = 398 },
    @{ Cell = "H26"; Value = 45261 },
    @{ Cell = "D27"; Value = 410 },
    @{ Cell = "H27"; Value = 45261 },
    @{ Cell = "C28"; Value = 359 },
    @{ Cell = "D28"; Value = 367 },
    @{ Cell = "F28"; Value = 45200 },
    @{ Cell = "H28"; Value = 45261 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
